$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Simple Fields")
$ws2 = $wb.Worksheets.Item("Simple Fields - Formatted")

# "Payment Terms" text tweak: "due 30 days" -> "30 days"
$ws1.Range("M2").Value = "30 days"
$ws2.Range("M2").Value = "30 days"

# Drop the "Currency" column (was column S / 19) from both Simple Fields sheets;
# this shifts the trailing "Items" column left from T into S.
$ws1.Range("S1:S2").EntireColumn.Delete()
$ws2.Range("S1:S2").EntireColumn.Delete()
